$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-12) got shuffled: each new row's D/J/K/L/M/O/P values
# are copied from a different original row. Mapping is new row -> source row
# (values taken from the row as it existed BEFORE any edits).
$data = @{
    2  = @{ D = 44650; J = 130; K = 3000; L = 3500; M = 3308; O = "Región Metropolitana";   P = 551 }
    3  = @{ D = 44637; J = 170; K = 2800; L = 3000; M = 2906; O = "Región Metropolitana";   P = 484 }
    4  = @{ D = 44659; J = 90;  K = 2500; L = 3000; M = 2722; O = "Región Metropolitana";   P = 454 }
    5  = @{ D = 44685; J = 150; K = 3000; L = 3500; M = 3267; O = "Región Metropolitana";   P = 544 }
    6  = @{ D = 44644; J = 140; K = 2500; L = 3000; M = 2786; O = "Provincia de Chacabuco"; P = 464 }
    7  = @{ D = 44630; J = 90;  K = 2500; L = 3000; M = 2722; O = "Región Metropolitana";   P = 454 }
    8  = @{ D = 44658; J = 180; K = 2500; L = 3000; M = 2778; O = "Región Metropolitana";   P = 463 }
    9  = @{ D = 44643; J = 90;  K = 2800; L = 3000; M = 2911; O = "Región Metropolitana";   P = 485 }
    10 = @{ D = 44631; J = 110; K = 3000; L = 3500; M = 3273; O = "Provincia de Chacabuco"; P = 546 }
    11 = @{ D = 44671; J = 150; K = 3500; L = 4000; M = 3733; O = "Región Metropolitana";   P = 622 }
    12 = @{ D = 44672; J = 140; K = 3000; L = 3500; M = 3286; O = "Región Metropolitana";   P = 548 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D    # D: Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K   # K: Precio mínimo
    $ws.Cells.Item($row, 12).Value = $vals.L   # L: Precio máximo
    $ws.Cells.Item($row, 13).Value = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $vals.O   # O: Origen
    $ws.Cells.Item($row, 16).Value = $vals.P   # P: Precio $/Kg
}
